$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of the Membrillo price series (Agro Chillan terminal).
# Full snapshot of rows 2-31 (A:T) as they should read after the refresh.
$data = @(
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45083, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 60, 9000, 10000, 9500, '$/caja 18 kilos empedrada', 'Región del Maule', 528, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45044, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 40, 13000, 13000, 13000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 722, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45044, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 40, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45043, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 40, 13000, 13000, 13000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 722, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45043, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 50, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44699, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 60, 13000, 13000, 13000, '$/caja 15 kilos granel', 'Provincia de Curicó', 867, 15),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44699, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 120, 11000, 12000, 11500, '$/caja 15 kilos granel', 'Provincia de Curicó', 767, 15),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45050, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 50, 13000, 13000, 13000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 722, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45050, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 40, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45071, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 40, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45071, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Segunda', 40, 10000, 10000, 10000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 556, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45049, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 50, 13000, 13000, 13000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 722, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45049, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 60, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45079, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 50, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45079, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 30, 10000, 10000, 10000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 556, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45079, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Segunda', 20, 9000, 9000, 9000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 500, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45069, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 60, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45069, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Segunda', 40, 10000, 10000, 10000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 556, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45020, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 60, 12000, 12000, 12000, '$/caja 18 kilos granel', 'Región de O''Higgins', 667, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45040, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 50, 13000, 13000, 13000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 722, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45040, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 40, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45070, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 60, 10000, 10000, 10000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 556, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45062, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 50, 13000, 13000, 13000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 722, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45062, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 50, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45085, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 50, 10000, 10000, 10000, '$/caja 18 kilos empedrada', 'Región del Maule', 556, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45033, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 60, 13000, 13000, 13000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 722, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45033, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 80, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45076, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 30, 12000, 12000, 12000, '$/caja 15 kilos granel', 'Región de O''Higgins', 800, 15),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45076, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Segunda', 30, 10000, 10000, 10000, '$/caja 15 kilos granel', 'Región de O''Higgins', 667, 15),
    @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45021, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 50, 12000, 12000, 12000, '$/caja 18 kilos granel', 'Región de O''Higgins', 667, 18)
)

$startRow = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($startRow, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $startRow = $startRow + 1
}

# Column D carries a date/time display format; make sure the newly created
# row (31) picks up the same number format as the rest of the date column.
$ws.Range("D31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
